$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Pipeline" debug signals added to the table (ack for debugger):
# 0x18 pipeline_immediate_out        (16 Bit)
# 0x19 pipeline_write_address_out    (4 Bit)
# 0x1A pipeline_whb_wlb_out          (2 Bit)
# 0x1D pipeline_is_alu_ram_gpu_op_out(3 Bit)

$ws.Range("B26").Value = "0x18"
$ws.Range("E26").Value = "pipeline_immediate_out"
$ws.Range("G26").Value = "0x18"
$ws.Range("H26").Value = "16 Bit"
$ws.Range("J26").Value = "pipeline_immediate_out"

$ws.Range("B27").Value = "0x19"
$ws.Range("E27").Value = "pipeline_write_address_out"
$ws.Range("G27").Value = "0x19"
$ws.Range("H27").Value = "4 Bit"
$ws.Range("J27").Value = "pipeline_write_address_out"

$ws.Range("B28").Value = "0x1A"
$ws.Range("E28").Value = "pipeline_whb_wlb_out"
$ws.Range("G28").Value = "0x1A"
$ws.Range("H28").Value = "2 Bit"
$ws.Range("J28").Value = "pipeline_whb_wlb_out"

$ws.Range("B31").Value = "0x1D"
$ws.Range("E31").Value = "pipeline_is_alu_ram_gpu_op_out"
$ws.Range("G31").Value = "0x1D"
$ws.Range("H31").Value = "3 Bit"
$ws.Range("J31").Value = "pipeline_is_alu_ram_gpu_op_out"

# Move the selection/active cell to B16 (also scrolls the view back to the top)
$ws.Range("B16").Select() | Out-Null
